$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (appended after existing row 93):
#   row 94 -> "dia" 93, date 30/06/2020
#   row 95 -> "dia" 94, date 01/07/2020
# Columns: A dia, B data, C casos, D mortes, E Ativos, F confirmado/100k,
#          G taxa morte, H Curados, I Casos negativos, J Testes realizados,
#          K novosCasos, L suspeitos, M mortesSuspeitas, N suspeitosAtivos,
#          O novosTestes, P leitos_clinicos_ocupados, Q leitos_uti_ocupados, R semana

$newRows = @(
    @{ Row = 94; A = 93; B = "30/06/2020"; C = 701; D = 19; E = 112; F = "461,524281";  G = "0,02710413695"; H = 567; I = 925; J = 1626; K = 35; L = 53; M = 5; N = 48; O = 94; P = 17; Q = 14; R = 14 },
    @{ Row = 95; A = 94; B = "01/07/2020"; C = 733; D = 21; E = 126; F = "482,5924365"; G = "0,02864938608"; H = 583; I = 989; J = 1722; K = 32; L = 46; M = 4; N = 42; O = 96; P = 22; Q = 14; R = 14 }
)

foreach ($entry in $newRows) {
    $r = $entry.Row
    $prevRow = $r - 1

    $ws.Cells.Item($r, 1).Value = $entry.A
    # Text-like columns: force literal text (avoid locale re-interpretation
    # of comma decimals / ambiguous dd/mm dates as numbers or dates).
    $ws.Cells.Item($r, 2).Value = "'" + $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D
    $ws.Cells.Item($r, 5).Value = $entry.E
    $ws.Cells.Item($r, 6).Value = "'" + $entry.F
    $ws.Cells.Item($r, 7).Value = "'" + $entry.G
    $ws.Cells.Item($r, 8).Value = $entry.H
    $ws.Cells.Item($r, 9).Value = $entry.I
    $ws.Cells.Item($r, 10).Value = $entry.J
    $ws.Cells.Item($r, 11).Value = $entry.K
    $ws.Cells.Item($r, 12).Value = $entry.L
    $ws.Cells.Item($r, 13).Value = $entry.M
    $ws.Cells.Item($r, 14).Value = $entry.N
    $ws.Cells.Item($r, 15).Value = $entry.O
    $ws.Cells.Item($r, 16).Value = $entry.P
    $ws.Cells.Item($r, 17).Value = $entry.Q
    $ws.Cells.Item($r, 18).Value = $entry.R

    # The apostrophe-prefixed entries above pick up a "quote prefix" style.
    # Restore plain (unstyled) formatting by copying the format from the
    # row above, matching the rest of the data rows which carry no style.
    $ws.Range($ws.Cells.Item($prevRow, 1), $ws.Cells.Item($prevRow, 18)).Copy()
    $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 18)).PasteSpecial(-4122)
}

$excel.CutCopyMode = $false
